# Solar_Data_Project_Briefing.docx edit:
#   "removed kw/h from required output and leaving in only dollar result"
#
# 1) The bullet "Savings expressed as kWh and dollars" loses the
#    "kWh and " portion, leaving "Savings expressed as dollars". Word
#    drops the floating "_GoBack" bookmark at the point of the last edit,
#    so we relocate it from its old (now-stale) spot to the new edit
#    point between the two remaining runs.
# 2) Incidental cleanup: the "(bill cost per kw/h)/60 for kw per minute"
#    sentence had its "kw" wrapped in proofing marks (w:proofErr) across
#    three runs; re-typing/normalizing that text merges it back into a
#    single clean run with no proofing marks, same wording.

$d = $word.ActiveDocument

# --- 1a. Remove the stale _GoBack bookmark wherever it currently sits ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 1b. Trim "kWh and " out of the savings bullet ---
$target = $d.Content
$target.Find.Execute("Savings expressed as kWh and dollars")
$lineStart = $target.Start

$prefix = "Savings expressed as "
$removed = "kWh and "

$cutStart = $lineStart + $prefix.Length
$cutEnd = $cutStart + $removed.Length

$cutRange = $d.Range($cutStart, $cutEnd)
$cutRange.Text = ""

# --- 1c. Re-plant _GoBack as a zero-length bookmark at the edit point ---
$editPoint = $d.Range($cutStart, $cutStart)
$d.Bookmarks.Add("_GoBack", $editPoint)

# --- 2. Normalize the "kw per minute" sentence back into one clean run ---
$d.Content.Find.Execute(" cost per kw/h)/60 for kw per minute", $false, $false, $false, $false, $false, $true, 1, $false, " cost per kw/h)/60 for kw per minute", 2)
